# New crime data collected — weekly CompStat update for the 123rd Precinct.
# Rolls the report forward one week (Volume/Number + reporting date range)
# and refreshes every crime-stat figure in the table (rows 16-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 32   Number  50" -> "...  51" ---------------------
$numberRun = $ws.Range("A8").Characters(21, 2)
$numberRun.Text = "51"

# --- Header: reporting week date range ----------------------------------
# "Report Covering the Week  12/8/2025  Through  12/14/2025"
# ->                         12/15/2025           12/21/2025
$weekStart = $ws.Range("C9").Characters(27, 9)
$weekStart.Text = "12/15/2025"
$weekEnd = $ws.Range("C9").Characters(48, 10)
$weekEnd.Text = "12/21/2025"

# --- Row 16 (Rape) -------------------------------------------------------
$ws.Range("C15").Copy($ws.Range("F16"))
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = -100
$ws.Range("L16").Value = 14.285714285714
$ws.Range("N16").Value = -62.5

# --- Row 17 (Robbery) -----------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -30
$ws.Range("I17").Value = 98
$ws.Range("J17").Value = 94
$ws.Range("K17").Value = 4.255319148936
$ws.Range("L17").Value = 15.294117647058
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -10.909090909090

# --- Row 18 (Fel. Assault) ------------------------------------------------
$ws.Range("D18").Value = 1
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = -75
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -21.052631578947
$ws.Range("N18").Value = -86.404833836858

# --- Row 19 (Burglary) -----------------------------------------------------
$ws.Range("C19").Value = 4
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 35.714285714285
$ws.Range("I19").Value = 283
$ws.Range("J19").Value = 256
$ws.Range("K19").Value = 10.546875
$ws.Range("L19").Value = -1.048951048951
$ws.Range("M19").Value = 93.835616438356
$ws.Range("N19").Value = 19.409282700421

# --- Row 20 (Gr. Larceny) ---------------------------------------------------
$ws.Range("I15").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 3
$ws.Range("I15").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 2
$ws.Range("K15").Copy($ws.Range("E20"))
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 29
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = -27.5
$ws.Range("L20").Value = -62.820512820512
$ws.Range("M20").Value = -23.684210526315
$ws.Range("N20").Value = -95.880681818181

# --- Row 21 (G.L.A.) ---------------------------------------------------------
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 32
$ws.Range("G21").Value = 34
$ws.Range("H21").Value = -5.882352941176
$ws.Range("I21").Value = 486
$ws.Range("J21").Value = 469
$ws.Range("K21").Value = 3.624733475479
$ws.Range("L21").Value = -8.301886792452
$ws.Range("M21").Value = 33.884297520661
$ws.Range("N21").Value = -66.505858028945

# --- Row 24 (Transit) ----------------------------------------------------
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = -72.727272727272
$ws.Range("F24").Value = 22
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = -46.341463414634
$ws.Range("I24").Value = 415
$ws.Range("J24").Value = 379
$ws.Range("K24").Value = 9.498680738786
$ws.Range("L24").Value = -17.165668662674
$ws.Range("M24").Value = -22.429906542056

# --- Row 25 (Housing) -----------------------------------------------------
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 13
$ws.Range("H25").Value = -23.529411764705
$ws.Range("I25").Value = 238
$ws.Range("J25").Value = 174
$ws.Range("K25").Value = 36.781609195402
$ws.Range("L25").Value = 17.241379310344

# --- Row 26 (Petit Larceny) -----------------------------------------------
$ws.Range("C15").Copy($ws.Range("C26"))
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 9
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = -30.769230769230
$ws.Range("J26").Value = 160
$ws.Range("K26").Value = 26.25
$ws.Range("M26").Value = -6.481481481481

# --- Row 28 (Misd. Assault) -----------------------------------------------
$ws.Range("C15").Copy($ws.Range("D28"))
$ws.Range("E15").Copy($ws.Range("E28"))
$ws.Range("L28").Value = -27.777777777777

# --- Row 29 (UCR Rape*) -----------------------------------------------------
$ws.Range("K15").Copy($ws.Range("M29"))
$ws.Range("M29").Value = -100

# --- Row 30 (Other Sex Crimes) ----------------------------------------------
$ws.Range("K15").Copy($ws.Range("M30"))
$ws.Range("M30").Value = -100

# --- Row 31 (Shooting Vic.) --------------------------------------------------
$ws.Range("C15").Copy($ws.Range("D31"))
$ws.Range("E15").Copy($ws.Range("E31"))

# --- Row 33 (Hate Crimes) -----------------------------------------------------
$ws.Range("C15").Copy($ws.Range("D33"))
$ws.Range("E15").Copy($ws.Range("E33"))
